$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reference cell with the default/normal style (style index 0), used to
# reset style on cells where a leading apostrophe was needed to keep
# numeric-looking text (e.g. "0.9995") stored as text instead of a number.
$normalStyle = $ws.Range("C2").Style

$ws.Range('D2').Value = '26.061.99'
$ws.Range('E2').Value = '  +0.81%  '
$ws.Range('D3').Value = '1.747.99'
$ws.Range('E3').Value = '  +0.29%  '
$ws.Range('D4').Value = "'0.9995"
$ws.Range('D4').Style = $normalStyle
$ws.Range('E4').Value = '  -0.09%  '
$ws.Range('D5').Value = "'234.89"
$ws.Range('D5').Style = $normalStyle
$ws.Range('E5').Value = '  +4.43%  '
$ws.Range('D6').Value = "'0.9995"
$ws.Range('D6').Style = $normalStyle
$ws.Range('E6').Value = '  -0.04%  '
$ws.Range('D7').Value = "'0.5287"
$ws.Range('D7').Style = $normalStyle
$ws.Range('E7').Value = '  +2.89%  '
$ws.Range('D8').Value = "'0.2786"
$ws.Range('D8').Style = $normalStyle
$ws.Range('E8').Value = '  +0.42%  '
$ws.Range('D9').Value = "'0.06183"
$ws.Range('D9').Style = $normalStyle
$ws.Range('E9').Value = '  +1.56%  '
$ws.Range('D10').Value = '1.747.25'
$ws.Range('E10').Value = '  +0.26%  '
$ws.Range('D11').Value = "'0.07180"
$ws.Range('D11').Style = $normalStyle
$ws.Range('E11').Value = '  +2.74%  '
$ws.Range('D12').Value = "'15.38"
$ws.Range('D12').Style = $normalStyle
$ws.Range('E12').Value = '  +1.08%  '
$ws.Range('D13').Value = "'0.6444"
$ws.Range('D13').Style = $normalStyle
$ws.Range('E13').Value = '  +1.56%  '
$ws.Range('D14').Value = "'4.621"
$ws.Range('D14').Style = $normalStyle
$ws.Range('E14').Value = '  +2.50%  '
$ws.Range('D15').Value = "'78.36"
$ws.Range('D15').Style = $normalStyle
$ws.Range('E15').Value = '  +2.40%  '
$ws.Range('D16').Value = "'0.9992"
$ws.Range('D16').Style = $normalStyle
$ws.Range('E16').Value = '  -0.05%  '
$ws.Range('D17').Value = "'0.9993"
$ws.Range('D17').Style = $normalStyle
$ws.Range('E17').Value = '  -0.10%  '
$ws.Range('D18').Value = '25.981.13'
$ws.Range('E18').Value = '  +0.42%  '
$ws.Range('D19').Value = "'11.69"
$ws.Range('D19').Style = $normalStyle
$ws.Range('E19').Value = '  +1.99%  '
$ws.Range('D20').Value = "'0.000006717"
$ws.Range('D20').Style = $normalStyle
$ws.Range('E20').Value = '  +2.07%  '
$ws.Range('D21').Value = '1.968.18'
$ws.Range('E21').Value = '  +0.52%  '
$ws.Range('D22').Value = "'4.311"
$ws.Range('D22').Style = $normalStyle
$ws.Range('E22').Value = '  +5.40%  '
$ws.Range('D23').Value = "'8.741"
$ws.Range('D23').Style = $normalStyle
$ws.Range('E23').Value = '  +2.75%  '
$ws.Range('D24').Value = "'5.220"
$ws.Range('D24').Style = $normalStyle
$ws.Range('D25').Value = "'138.12"
$ws.Range('D25').Style = $normalStyle
$ws.Range('D26').Value = "'1.508"
$ws.Range('D26').Style = $normalStyle
$ws.Range('E26').Value = '  +0.54%  '
$ws.Range('E27').Value = '  +2.16%  '
$ws.Range('D28').Value = "'1.804"
$ws.Range('D28').Style = $normalStyle
$ws.Range('E28').Value = '  -0.88%  '
$ws.Range('D29').Value = "'104.48"
$ws.Range('D29').Style = $normalStyle
$ws.Range('E29').Value = '  +1.80%  '
$ws.Range('D30').Value = "'0.08276"
$ws.Range('D30').Style = $normalStyle
$ws.Range('E30').Value = '  +0.00%  '
$ws.Range('D31').Value = "'3.806"
$ws.Range('D31').Style = $normalStyle
$ws.Range('E31').Value = '  +5.06%  '
$ws.Range('D32').Value = "'3.660"
$ws.Range('D32').Style = $normalStyle
$ws.Range('E32').Value = '  +7.69%  '
$ws.Range('D33').Value = "'0.04574"
$ws.Range('D33').Style = $normalStyle
$ws.Range('E33').Value = '  +3.86%  '
$ws.Range('D34').Value = "'2.642"
$ws.Range('D34').Style = $normalStyle
$ws.Range('E34').Value = '  +0.93%  '
$ws.Range('E35').Value = '  +3.45%  '
$ws.Range('D36').Value = "'0.6346"
$ws.Range('D36').Style = $normalStyle
$ws.Range('E36').Value = '  +6.11%  '
$ws.Range('D37').Value = "'2.714"
$ws.Range('D37').Style = $normalStyle
$ws.Range('E37').Value = '  +1.54%  '
$ws.Range('E38').Value = '  +2.73%  '
$ws.Range('D39').Value = "'1.942"
$ws.Range('D39').Style = $normalStyle
$ws.Range('E39').Value = '  +1.59%  '
$ws.Range('D40').Value = "'0.9989"
$ws.Range('D40').Style = $normalStyle
$ws.Range('E40').Value = '  +0.01%  '
$ws.Range('D41').Value = "'100.13"
$ws.Range('D41').Style = $normalStyle
$ws.Range('E41').Value = '  -0.76%  '
$ws.Range('D42').Value = "'0.3923"
$ws.Range('D42').Style = $normalStyle
$ws.Range('E42').Value = '  +2.39%  '
$ws.Range('D43').Value = "'0.7458"
$ws.Range('D43').Style = $normalStyle
$ws.Range('E43').Value = '  +3.03%  '
$ws.Range('D44').Value = "'5.031"
$ws.Range('D44').Style = $normalStyle
$ws.Range('E44').Value = '  +3.28%  '
$ws.Range('D45').Value = "'0.1145"
$ws.Range('D45').Style = $normalStyle
$ws.Range('E45').Value = '  +3.92%  '
$ws.Range('D46').Value = "'6.322"
$ws.Range('D46').Style = $normalStyle
$ws.Range('E46').Value = '  +1.09%  '
$ws.Range('D47').Value = "'0.05346"
$ws.Range('D47').Style = $normalStyle
$ws.Range('E47').Value = '  -2.27%  '
$ws.Range('D48').Value = "'30.86"
$ws.Range('D48').Style = $normalStyle
$ws.Range('E48').Value = '  +3.95%  '
$ws.Range('D49').Value = "'54.11"
$ws.Range('D49').Style = $normalStyle
$ws.Range('E49').Value = '  +3.61%  '
$ws.Range('D50').Value = "'7.638"
$ws.Range('D50').Style = $normalStyle
$ws.Range('E50').Value = '  +2.31%  '
$ws.Range('D51').Value = "'0.3456"
$ws.Range('D51').Style = $normalStyle
$ws.Range('E51').Value = '  +1.91%  '
